$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Running final trends estimates: refine theta/lambda/proportion-drinking
# figures to higher precision.

# Row 2 (1983)
$ws.Range("E2").Value = 0.182455

# Row 3 (1987)
$ws.Range("D3").Value = "(0.03)"
$ws.Range("E3").Value = "(0.00001)"

# Row 4 (1988)
$ws.Range("E4").Value = 0.167269

# Row 5 (1992)
$ws.Range("D5").Value = "(0.05)"
$ws.Range("E5").Value = "(0.00001)"

# Row 6 (1993)
$ws.Range("E6").Value = 0.122975

# Row 7 (1997)
$ws.Range("C7").Value = "(0.15)"
$ws.Range("D7").Value = "(0.09)"
$ws.Range("E7").Value = "(0.00001)"

# Row 8 (1998)
$ws.Range("E8").Value = 0.101189

# Row 9 (2002)
$ws.Range("C9").Value = "(0.25)"
$ws.Range("D9").Value = "(0.17)"
$ws.Range("E9").Value = "(0.00001)"

# Row 10 (2003)
$ws.Range("E10").Value = 0.092802

# Row 11 (2007)
$ws.Range("C11").Value = "(0.25)"
$ws.Range("D11").Value = "(0.22)"
$ws.Range("E11").Value = "(0.00001)"

# Row 12 (2008)
$ws.Range("E12").Value = 0.097228

# Row 13 (2012)
$ws.Range("C13").Value = "(0.31)"
$ws.Range("D13").Value = "(0.21)"
$ws.Range("E13").Value = "(0.00001)"

# Row 14 (2013)
$ws.Range("E14").Value = 0.089807

# Row 15 (2017)
$ws.Range("C15").Value = "(0.31)"
$ws.Range("D15").Value = "(0.17)"
$ws.Range("E15").Value = "(0.00001)"
